$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell for new column K, matching style of existing header row (bold/border/center)
$ws.Range("K1").Value = "intervention_type"
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# Rows where intervention_type is present-but-empty in the source data (an
# empty inline-string cell, rather than the column being entirely absent).
# Paste an existing empty inline-string cell's format (D2, which is blank in
# the source) onto each of these so the K cell exists but stays empty.
$ws.Range("D2").Copy()
foreach ($r in @(58, 99, 100, 101, 102, 103)) {
    $ws.Range("K" + $r).PasteSpecial(-4122)
}

# Data for intervention_type column (row => value).
$data = @(
    @{Row=2; Value="DRUG"},
    @{Row=3; Value="DEVICE"},
    @{Row=4; Value="PROCEDURE"},
    @{Row=5; Value="DRUG"},
    @{Row=6; Value="PROCEDURE"},
    @{Row=7; Value="DRUG"},
    @{Row=8; Value="PROCEDURE"},
    @{Row=9; Value="DRUG"},
    @{Row=10; Value="PROCEDURE"},
    @{Row=11; Value="DEVICE"},
    @{Row=12; Value="DRUG"},
    @{Row=13; Value="PROCEDURE"},
    @{Row=14; Value="OTHER"},
    @{Row=15; Value="DRUG"},
    @{Row=16; Value="DRUG"},
    @{Row=17; Value="PROCEDURE"},
    @{Row=18; Value="PROCEDURE"},
    @{Row=19; Value="DEVICE"},
    @{Row=20; Value="PROCEDURE"},
    @{Row=21; Value="OTHER"},
    @{Row=22; Value="OTHER"},
    @{Row=23; Value="PROCEDURE"},
    @{Row=24; Value="DIETARY_SUPPLEMENT"},
    @{Row=25; Value="DEVICE"},
    @{Row=26; Value="DEVICE"},
    @{Row=27; Value="GENETIC"},
    @{Row=28; Value="PROCEDURE"},
    @{Row=29; Value="PROCEDURE"},
    @{Row=30; Value="OTHER"},
    @{Row=31; Value="DEVICE"},
    @{Row=32; Value="OTHER"},
    @{Row=33; Value="DRUG"},
    @{Row=34; Value="DEVICE"},
    @{Row=35; Value="BEHAVIORAL"},
    @{Row=36; Value="DRUG"},
    @{Row=37; Value="DEVICE"},
    @{Row=38; Value="DEVICE"},
    @{Row=39; Value="PROCEDURE"},
    @{Row=40; Value="RADIATION"},
    @{Row=41; Value="DEVICE"},
    @{Row=42; Value="DRUG"},
    @{Row=43; Value="BEHAVIORAL"},
    @{Row=44; Value="BEHAVIORAL"},
    @{Row=45; Value="PROCEDURE"},
    @{Row=46; Value="OTHER"},
    @{Row=47; Value="OTHER"},
    @{Row=48; Value="PROCEDURE"},
    @{Row=49; Value="OTHER"},
    @{Row=50; Value="DIAGNOSTIC_TEST"},
    @{Row=51; Value="OTHER"},
    @{Row=52; Value="BIOLOGICAL"},
    @{Row=53; Value="DEVICE"},
    @{Row=54; Value="DIAGNOSTIC_TEST"},
    @{Row=55; Value="OTHER"},
    @{Row=56; Value="OTHER"},
    @{Row=57; Value="DIAGNOSTIC_TEST"},
    @{Row=59; Value="DIAGNOSTIC_TEST"},
    @{Row=60; Value="OTHER"},
    @{Row=61; Value="DRUG"},
    @{Row=62; Value="DIAGNOSTIC_TEST"},
    @{Row=63; Value="BIOLOGICAL"},
    @{Row=64; Value="DEVICE"},
    @{Row=65; Value="OTHER"},
    @{Row=66; Value="DEVICE"},
    @{Row=67; Value="DRUG"},
    @{Row=68; Value="OTHER"},
    @{Row=69; Value="DEVICE"},
    @{Row=70; Value="OTHER"},
    @{Row=71; Value="OTHER"},
    @{Row=72; Value="OTHER"},
    @{Row=73; Value="DEVICE"},
    @{Row=74; Value="OTHER"},
    @{Row=75; Value="DEVICE"},
    @{Row=76; Value="OTHER"},
    @{Row=77; Value="OTHER"},
    @{Row=78; Value="DIAGNOSTIC_TEST"},
    @{Row=79; Value="OTHER"},
    @{Row=80; Value="PROCEDURE"},
    @{Row=81; Value="OTHER"},
    @{Row=82; Value="DIAGNOSTIC_TEST"},
    @{Row=83; Value="OTHER"},
    @{Row=84; Value="OTHER"},
    @{Row=85; Value="OTHER"},
    @{Row=86; Value="OTHER"},
    @{Row=87; Value="PROCEDURE"},
    @{Row=88; Value="DIAGNOSTIC_TEST"},
    @{Row=89; Value="DIAGNOSTIC_TEST"},
    @{Row=90; Value="DRUG"},
    @{Row=91; Value="OTHER"},
    @{Row=92; Value="DIETARY_SUPPLEMENT"},
    @{Row=93; Value="PROCEDURE"},
    @{Row=94; Value="DIAGNOSTIC_TEST"},
    @{Row=95; Value="PROCEDURE"},
    @{Row=96; Value="PROCEDURE"},
    @{Row=97; Value="DRUG"},
    @{Row=98; Value="DRUG"}
)

foreach ($entry in $data) {
    $ws.Cells.Item($entry.Row, 11).Value = $entry.Value
}
